$wb = $excel.ActiveWorkbook

# --- Update DSPChannelMap data values (Channel_Number_DSP2 / B column renumbering) ---
$ws = $wb.Worksheets.Item("DSPChannelMap")

$ws.Range("D5").Value = 10
$ws.Range("D6").Value = 11
$ws.Range("D7").Value = 12
$ws.Range("D8").Value = 13
$ws.Range("D9").Value = 14
$ws.Range("D10").Value = 15

$ws.Range("B11").Value = 16
$ws.Range("B12").Value = 17
$ws.Range("B13").Value = 18

# --- Update sheet selections / active tab ---
# DeviceInfo: clear its tabSelected state and move the selection to I14
$wsDeviceInfo = $wb.Worksheets.Item("DeviceInfo")
$wsDeviceInfo.Activate()
$wsDeviceInfo.Range("I14").Select()

# DSPChannelMap becomes the active tab, with B13 selected
$ws.Activate()
$ws.Range("B13").Select()
